$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.414.98"
$ws.Range("E2").Value = "  +0.00%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.708.73"
$ws.Range("E3").Value = "  +2.26%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.87"
$ws.Range("E5").Value = "  +1.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.77"
$ws.Range("E6").Value = "  +4.95%  "

# Row 7
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.558"
$ws.Range("E8").Value = "  +3.38%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.709.03"
$ws.Range("E9").Value = "  +2.30%  "

# Row 10
$ws.Range("E10").Value = "  +2.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  +0.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.363"
$ws.Range("E12").Value = "  +3.36%  "

# Row 13
$ws.Range("E13").Value = "  +0.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.39"
$ws.Range("E14").Value = "  +1.36%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.204.21"
$ws.Range("E15").Value = "  +2.28%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000187"
$ws.Range("E16").Value = "  -0.16%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.326.09"
$ws.Range("E17").Value = "  +0.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.713.87"
$ws.Range("E18").Value = "  +2.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.79"
$ws.Range("E19").Value = "  +2.12%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "370.98"
$ws.Range("E20").Value = "  +1.99%  "

# Row 21
$ws.Range("E21").Value = "  +1.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.48"
$ws.Range("E22").Value = "  +1.62%  "

# Row 23
$ws.Range("E23").Value = "  +3.37%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -0.50%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.17"
$ws.Range("E25").Value = "  -2.08%  "

# Row 26
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.00"
$ws.Range("E27").Value = "  +1.08%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.859.79"
$ws.Range("E28").Value = "  +2.73%  "

# Row 29
$ws.Range("E29").Value = "  +0.57%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "579.95"
$ws.Range("E31").Value = "  +2.74%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.10"
$ws.Range("E32").Value = "  +0.43%  "

# Row 33
$ws.Range("E33").Value = "  +0.94%  "

# Row 34
$ws.Range("E34").Value = "  +5.40%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.131"

# Row 36
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.59"
$ws.Range("E37").Value = "  -3.39%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.26"
$ws.Range("E38").Value = "  +0.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.85"
$ws.Range("E39").Value = "  +1.03%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.377"
$ws.Range("E40").Value = "  +1.84%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.87"
$ws.Range("E41").Value = "  -0.07%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("E42").Value = "  +0.76%  "

# Row 43
$ws.Range("E43").Value = "  +1.00%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").Value = "  -1.67%  "

# Row 45
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0310"
$ws.Range("E46").Value = "  -3.26%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.79"
$ws.Range("E47").Value = "  +1.15%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.596"
$ws.Range("E48").Value = "  +3.60%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "154.36"
$ws.Range("E49").Value = "  -2.43%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.89"
$ws.Range("E50").Value = "  +1.77%  "

# Row 51
$ws.Range("E51").Value = "  +4.01%  "
